$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) and Volume(1h) (E) columns to be treated as plain text so
# that values such as "1.00", "0.200" or "0.0000200" keep their exact, literal
# formatting instead of being auto-converted into numbers (which would drop
# meaningful trailing/leading zeros).
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "99.376.09"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "3.292.22"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "253.38"
$ws.Range("E5").Value = "  -0.94%  "
$ws.Range("D6").Value = "623.22"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "1.42"
$ws.Range("E7").Value = "  +18.08%  "
$ws.Range("D8").Value = "0.399"
$ws.Range("E8").Value = "  +3.54%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "0.963"
$ws.Range("E10").Value = "  +20.05%  "
$ws.Range("D11").Value = "3.286.43"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").Value = "0.200"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Value = "39.48"
$ws.Range("E13").Value = "  +10.34%  "
$ws.Range("D14").Value = "99.084.60"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "0.0000247"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("D16").Value = "3.875.74"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "5.46"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").Value = "3.280.77"
$ws.Range("E18").Value = "  -2.19%  "
$ws.Range("D19").Value = "3.43"
$ws.Range("E19").Value = "  -3.66%  "
$ws.Range("D20").Value = "15.37"
$ws.Range("E20").Value = "  +2.85%  "
$ws.Range("D21").Value = "6.33"
$ws.Range("E21").Value = "  +8.51%  "
$ws.Range("D22").Value = "488.81"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").Value = "9.32"
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("D24").Value = "0.0000200"
$ws.Range("E24").Value = "  -2.54%  "
$ws.Range("D25").Value = "5.64"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("D26").Value = "88.78"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("D27").Value = "0.320"
$ws.Range("E27").Value = "  +27.90%  "
$ws.Range("D28").Value = "11.99"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").Value = "3.424.29"
$ws.Range("E29").Value = "  -2.98%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("D31").Value = "0.190"
$ws.Range("E31").Value = "  +1.25%  "
$ws.Range("D32").Value = "0.136"
$ws.Range("E32").Value = "  +9.02%  "
$ws.Range("D33").Value = "10.36"
$ws.Range("E33").Value = "  +11.88%  "
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "27.90"
$ws.Range("E35").Value = "  +2.53%  "
$ws.Range("D36").Value = "0.149"
$ws.Range("E36").Value = "  -1.21%  "
$ws.Range("D37").Value = "0.473"
$ws.Range("E37").Value = "  +6.44%  "
$ws.Range("D38").Value = "7.21"
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("D39").Value = "1.94"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").Value = "489.52"
$ws.Range("E41").Value = "  -5.22%  "
$ws.Range("D42").Value = "3.63"
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("D43").Value = "1.23"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "0.772"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "3.09"
$ws.Range("E46").Value = "  -4.55%  "
$ws.Range("D47").Value = "1.95"
$ws.Range("E47").Value = "  +2.46%  "
$ws.Range("D48").Value = "157.72"
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("D49").Value = "0.844"
$ws.Range("E49").Value = "  +6.99%  "
$ws.Range("D50").Value = "7.27"
$ws.Range("E50").Value = "  +15.76%  "
$ws.Range("D51").Value = "4.72"
$ws.Range("E51").Value = "  +4.93%  "
